$wb = $excel.ActiveWorkbook

# Cell value updates: (Sheet, CellRef, NewValue)
$updates = @(
    @{Sheet="ALC"; Cell="H43"; Value=1198},
    @{Sheet="ALC"; Cell="I43"; Value=1000},
    @{Sheet="ALC"; Cell="J43"; Value=1330},
    @{Sheet="ALC"; Cell="K43"; Value=1000},
    @{Sheet="ALC"; Cell="L43"; Value=1330},
    @{Sheet="ALC"; Cell="M43"; Value=-931},
    @{Sheet="ALC"; Cell="N43"; Value=-1468},
    @{Sheet="ALC"; Cell="H138"; Value=2056.53},
    @{Sheet="ALC"; Cell="I138"; Value=984.0909},
    @{Sheet="ALC"; Cell="J138"; Value=2189.0786},
    @{Sheet="ALC"; Cell="K138"; Value=2952.2727},
    @{Sheet="ALC"; Cell="L138"; Value=6567.235799999999},
    @{Sheet="ALC"; Cell="M138"; Value=2187.7273},
    @{Sheet="ALC"; Cell="N138"; Value=-16847.2358},
    @{Sheet="ARM"; Cell="H61"; Value=2384.5},
    @{Sheet="ARM"; Cell="I61"; Value=0},
    @{Sheet="ARM"; Cell="J61"; Value=2384.5},
    @{Sheet="ARM"; Cell="K61"; Value=0},
    @{Sheet="ARM"; Cell="L61"; Value=2384.5},
    @{Sheet="ARM"; Cell="N61"; Value=-2808.5},
    @{Sheet="ARM"; Cell="H74"; Value=596.71155},
    @{Sheet="ARM"; Cell="I74"; Value=566.61536},
    @{Sheet="ARM"; Cell="J74"; Value=687},
    @{Sheet="ARM"; Cell="K74"; Value=566.61536},
    @{Sheet="ARM"; Cell="L74"; Value=687},
    @{Sheet="ARM"; Cell="M74"; Value=307.38464},
    @{Sheet="ARM"; Cell="N74"; Value=-2435},
    @{Sheet="ARM"; Cell="H77"; Value=596.71155},
    @{Sheet="ARM"; Cell="I77"; Value=566.61536},
    @{Sheet="ARM"; Cell="J77"; Value=687},
    @{Sheet="ARM"; Cell="K77"; Value=2833.0768},
    @{Sheet="ARM"; Cell="L77"; Value=3435},
    @{Sheet="ARM"; Cell="M77"; Value=1534.9232},
    @{Sheet="ARM"; Cell="N77"; Value=-12171},
    @{Sheet="ARM"; Cell="H106"; Value=21185},
    @{Sheet="ARM"; Cell="J106"; Value=21185},
    @{Sheet="ARM"; Cell="L106"; Value=21185},
    @{Sheet="ARM"; Cell="N106"; Value=-23709},
    @{Sheet="ARM"; Cell="H132"; Value=4273.5713},
    @{Sheet="ARM"; Cell="I132"; Value=5045.4287},
    @{Sheet="ARM"; Cell="J132"; Value=3501.7144},
    @{Sheet="ARM"; Cell="K132"; Value=15136.2861},
    @{Sheet="ARM"; Cell="L132"; Value=10505.1432},
    @{Sheet="ARM"; Cell="M132"; Value=-12606.2861},
    @{Sheet="ARM"; Cell="N132"; Value=-15565.1432},
    @{Sheet="ARM"; Cell="H136"; Value=2384.5},
    @{Sheet="ARM"; Cell="I136"; Value=0},
    @{Sheet="ARM"; Cell="J136"; Value=2384.5},
    @{Sheet="ARM"; Cell="K136"; Value=0},
    @{Sheet="ARM"; Cell="L136"; Value=7153.5},
    @{Sheet="ARM"; Cell="N136"; Value=-12253.5},
    @{Sheet="BSM"; Cell="H31"; Value=500},
    @{Sheet="BSM"; Cell="I31"; Value=500},
    @{Sheet="BSM"; Cell="K31"; Value=500},
    @{Sheet="BSM"; Cell="M31"; Value=-248},
    @{Sheet="BSM"; Cell="H98"; Value=64333.332},
    @{Sheet="BSM"; Cell="J98"; Value=75000},
    @{Sheet="BSM"; Cell="L98"; Value=75000},
    @{Sheet="BSM"; Cell="N98"; Value=-80990},
    @{Sheet="BSM"; Cell="H99"; Value=26317090},
    @{Sheet="BSM"; Cell="I99"; Value=33334574},
    @{Sheet="BSM"; Cell="K99"; Value=33334574},
    @{Sheet="BSM"; Cell="M99"; Value=-33333076},
    @{Sheet="BSM"; Cell="H113"; Value=16666666},
    @{Sheet="BSM"; Cell="I113"; Value=16666666},
    @{Sheet="BSM"; Cell="K113"; Value=16666666},
    @{Sheet="BSM"; Cell="M113"; Value=-16664496},
    @{Sheet="BSM"; Cell="H134"; Value=10619.174},
    @{Sheet="BSM"; Cell="I134"; Value=8015.125},
    @{Sheet="BSM"; Cell="K134"; Value=24045.375},
    @{Sheet="BSM"; Cell="M134"; Value=-21510.375},
    @{Sheet="CRP"; Cell="H99"; Value=2203.5715},
    @{Sheet="CRP"; Cell="I99"; Value=2142},
    @{Sheet="CRP"; Cell="J99"; Value=2285.6667},
    @{Sheet="CRP"; Cell="K99"; Value=2142},
    @{Sheet="CRP"; Cell="L99"; Value=2285.6667},
    @{Sheet="CRP"; Cell="M99"; Value=-644},
    @{Sheet="CRP"; Cell="N99"; Value=-5281.6667},
    @{Sheet="CRP"; Cell="H105"; Value=1011},
    @{Sheet="CRP"; Cell="I105"; Value=0},
    @{Sheet="CRP"; Cell="J105"; Value=1011},
    @{Sheet="CRP"; Cell="K105"; Value=0},
    @{Sheet="CRP"; Cell="L105"; Value=1011},
    @{Sheet="CRP"; Cell="N105"; Value=-4505},
    @{Sheet="CRP"; Cell="H106"; Value=0},
    @{Sheet="CRP"; Cell="J106"; Value=0},
    @{Sheet="CRP"; Cell="L106"; Value=0},
    @{Sheet="CRP"; Cell="H126"; Value=2203.5715},
    @{Sheet="CRP"; Cell="I126"; Value=2142},
    @{Sheet="CRP"; Cell="J126"; Value=2285.6667},
    @{Sheet="CRP"; Cell="K126"; Value=6426},
    @{Sheet="CRP"; Cell="L126"; Value=6857.000100000001},
    @{Sheet="CRP"; Cell="M126"; Value=-3956},
    @{Sheet="CRP"; Cell="N126"; Value=-11797.0001},
    @{Sheet="CRP"; Cell="H131"; Value=30000},
    @{Sheet="CRP"; Cell="J131"; Value=30000},
    @{Sheet="CRP"; Cell="L131"; Value=30000},
    @{Sheet="CRP"; Cell="N131"; Value=-40080},
    @{Sheet="CUL"; Cell="H4"; Value=81352.25999999999},
    @{Sheet="CUL"; Cell="I4"; Value=166.625},
    @{Sheet="CUL"; Cell="J4"; Value=167950.27},
    @{Sheet="CUL"; Cell="K4"; Value=499.875},
    @{Sheet="CUL"; Cell="L4"; Value=503850.8099999999},
    @{Sheet="CUL"; Cell="M4"; Value=-387.875},
    @{Sheet="CUL"; Cell="N4"; Value=-504074.8099999999},
    @{Sheet="GSM"; Cell="H99"; Value=5250},
    @{Sheet="GSM"; Cell="I99"; Value=5250},
    @{Sheet="GSM"; Cell="K99"; Value=5250},
    @{Sheet="GSM"; Cell="M99"; Value=-3004},
    @{Sheet="GSM"; Cell="H100"; Value=40000},
    @{Sheet="GSM"; Cell="J100"; Value=40000},
    @{Sheet="GSM"; Cell="L100"; Value=40000},
    @{Sheet="GSM"; Cell="N100"; Value=-42164},
    @{Sheet="GSM"; Cell="H101"; Value=16328.5},
    @{Sheet="GSM"; Cell="J101"; Value=16328.5},
    @{Sheet="GSM"; Cell="L101"; Value=16328.5},
    @{Sheet="GSM"; Cell="N101"; Value=-22818.5},
    @{Sheet="GSM"; Cell="H107"; Value=756.0625},
    @{Sheet="GSM"; Cell="I107"; Value=643.7},
    @{Sheet="GSM"; Cell="K107"; Value=643.7},
    @{Sheet="GSM"; Cell="M107"; Value=1276.3},
    @{Sheet="GSM"; Cell="H122"; Value=2267.353},
    @{Sheet="GSM"; Cell="I122"; Value=1364.7},
    @{Sheet="GSM"; Cell="K122"; Value=4094.1},
    @{Sheet="GSM"; Cell="M122"; Value=-1644.1},
    @{Sheet="GSM"; Cell="H132"; Value=2810.5652},
    @{Sheet="GSM"; Cell="I132"; Value=2557.5454},
    @{Sheet="GSM"; Cell="K132"; Value=7672.6362},
    @{Sheet="GSM"; Cell="M132"; Value=-5142.6362},
    @{Sheet="LTW"; Cell="H7"; Value=2170.5},
    @{Sheet="LTW"; Cell="I7"; Value=2198.8333},
    @{Sheet="LTW"; Cell="K7"; Value=2198.8333},
    @{Sheet="LTW"; Cell="M7"; Value=-2086.8333},
    @{Sheet="LTW"; Cell="H22"; Value=1599.9166},
    @{Sheet="LTW"; Cell="I22"; Value=1425},
    @{Sheet="LTW"; Cell="J22"; Value=1949.75},
    @{Sheet="LTW"; Cell="K22"; Value=1425},
    @{Sheet="LTW"; Cell="L22"; Value=1949.75},
    @{Sheet="LTW"; Cell="M22"; Value=-1130},
    @{Sheet="LTW"; Cell="N22"; Value=-2539.75},
    @{Sheet="LTW"; Cell="H27"; Value=1599.9166},
    @{Sheet="LTW"; Cell="I27"; Value=1425},
    @{Sheet="LTW"; Cell="J27"; Value=1949.75},
    @{Sheet="LTW"; Cell="K27"; Value=1425},
    @{Sheet="LTW"; Cell="L27"; Value=1949.75},
    @{Sheet="LTW"; Cell="M27"; Value=-1318},
    @{Sheet="LTW"; Cell="N27"; Value=-2163.75},
    @{Sheet="LTW"; Cell="H40"; Value=2622.5},
    @{Sheet="LTW"; Cell="I40"; Value=2630},
    @{Sheet="LTW"; Cell="J40"; Value=2600},
    @{Sheet="LTW"; Cell="K40"; Value=2630},
    @{Sheet="LTW"; Cell="L40"; Value=2600},
    @{Sheet="LTW"; Cell="M40"; Value=-2494},
    @{Sheet="LTW"; Cell="N40"; Value=-2872},
    @{Sheet="LTW"; Cell="H123"; Value=37714.5},
    @{Sheet="LTW"; Cell="J123"; Value=37714.5},
    @{Sheet="LTW"; Cell="L123"; Value=37714.5},
    @{Sheet="LTW"; Cell="N123"; Value=-47514.5},
    @{Sheet="LTW"; Cell="H126"; Value=2170.5},
    @{Sheet="LTW"; Cell="I126"; Value=2198.8333},
    @{Sheet="LTW"; Cell="K126"; Value=6596.499899999999},
    @{Sheet="LTW"; Cell="M126"; Value=-4126.499899999999},
    @{Sheet="WVR"; Cell="H81"; Value=2180},
    @{Sheet="WVR"; Cell="I81"; Value=833.3333},
    @{Sheet="WVR"; Cell="K81"; Value=1666.6666},
    @{Sheet="WVR"; Cell="M81"; Value=-605.6666},
    @{Sheet="WVR"; Cell="H84"; Value=2180},
    @{Sheet="WVR"; Cell="I84"; Value=833.3333},
    @{Sheet="WVR"; Cell="K84"; Value=8333.333000000001},
    @{Sheet="WVR"; Cell="M84"; Value=-3029.333000000001},
    @{Sheet="WVR"; Cell="H96"; Value=1829.2727},
    @{Sheet="WVR"; Cell="I96"; Value=1682.4445},
    @{Sheet="WVR"; Cell="J96"; Value=2490},
    @{Sheet="WVR"; Cell="K96"; Value=1682.4445},
    @{Sheet="WVR"; Cell="L96"; Value=2490},
    @{Sheet="WVR"; Cell="M96"; Value=-309.4445000000001},
    @{Sheet="WVR"; Cell="N96"; Value=-5236},
    @{Sheet="WVR"; Cell="H122"; Value=7648170.5},
    @{Sheet="WVR"; Cell="I122"; Value=8667650},
    @{Sheet="WVR"; Cell="K122"; Value=26002950},
    @{Sheet="WVR"; Cell="M122"; Value=-26000500},
    @{Sheet="WVR"; Cell="H123"; Value=75000},
    @{Sheet="WVR"; Cell="J123"; Value=75000},
    @{Sheet="WVR"; Cell="L123"; Value=75000},
    @{Sheet="WVR"; Cell="N123"; Value=-84800},
    @{Sheet="WVR"; Cell="H126"; Value=47620748},
    @{Sheet="WVR"; Cell="I126"; Value=142858300},
    @{Sheet="WVR"; Cell="J126"; Value=1971.7858},
    @{Sheet="WVR"; Cell="K126"; Value=428574900},
    @{Sheet="WVR"; Cell="L126"; Value=5915.357400000001},
    @{Sheet="WVR"; Cell="M126"; Value=-428572430},
    @{Sheet="WVR"; Cell="N126"; Value=-10855.3574},
    @{Sheet="WVR"; Cell="H127"; Value=82803},
    @{Sheet="WVR"; Cell="J127"; Value=82803},
    @{Sheet="WVR"; Cell="L127"; Value=82803},
    @{Sheet="WVR"; Cell="N127"; Value=-92723}
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}

# Cells that must be cleared entirely (removed from the sheet)
$clears = @(
    @{Sheet="ARM"; Cell="M61"},
    @{Sheet="ARM"; Cell="M136"},
    @{Sheet="CRP"; Cell="M105"},
    @{Sheet="CRP"; Cell="N106"}
)

foreach ($c in $clears) {
    $ws = $wb.Worksheets.Item($c.Sheet)
    $ws.Range($c.Cell).ClearContents()
}

Write-Output "Applied $($updates.Count) updates and $($clears.Count) clears."
